$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 14
$ws.Range("B9").Value = "Nguyễn Thị B"
$ws.Range("C9").Value = "Hà Nội"
$ws.Range("D9").Value = "'0123456789"
$ws.Range("D9").ClearFormats()

# Row 10
$ws.Range("A10").Value = 15
$ws.Range("B10").Value = "'1"
$ws.Range("C10").Value = "'1"
$ws.Range("D10").Value = "'1"
$ws.Range("B10:D10").ClearFormats()

# Row 11
$ws.Range("A11").Value = 16
$ws.Range("B11").Value = "'2"
$ws.Range("C11").Value = "'2"
$ws.Range("D11").Value = "'2"
$ws.Range("B11:D11").ClearFormats()
